$d = $word.ActiveDocument
$enDash = [char]0x2013

# Replace the date-range string (appears multiple times, identical in each case)
$d.Content.Find.Execute(
    "2022: Datumi kampanje za opazovanje Ozvezdje Laboda: 10.-19. avgust, 9.-18. september, 8.-17. oktober",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ("2022: Datumi kampanje za opazovanje ozvezdje Cygnus: 10." + $enDash + "19. avgust, 9." + $enDash + "18. september, 8." + $enDash + "17. oktober"),
    2
)

# Replace the inline mention within the descriptive paragraph
$d.Content.Find.Execute(
    "izbranega Ozvezdje Laboda na nočnem nebu",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "izbranega ozvezdje Cygnus na nočnem nebu",
    2
)
